# Lab01 review report update - "Updated Lab01 once againnn"
$wb = $excel.ActiveWorkbook

# --- Sheet: Requirements Phase Defects ---
$ws1 = $wb.Worksheets.Item("Requirements Phase Defects")
$ws1.Range("C14").Value = $null
$ws1.Range("E14").Value = $null
[void]$ws1.Range("C14").Select()

# --- Sheet: Architect. Design Phase Defects ---
$ws2 = $wb.Worksheets.Item("Architect. Design Phase Defects")
$ws2.Range("C12").Value = $null
$ws2.Range("E12").Value = $null
$ws2.Range("C13").Value = $null
$ws2.Range("E13").Value = $null
$ws2.Range("C15").Value = $null
$ws2.Range("E15").Value = $null
$ws2.Range("C17").Value = $null
$ws2.Range("E17").Value = $null
$ws2.Range("C18").Value = $null
$ws2.Range("E18").Value = $null
[void]$ws2.Range("C18").Select()

# --- Sheet: Coding Phase Defects ---
$ws3 = $wb.Worksheets.Item("Coding Phase Defects")
$ws3.Range("C11").Value = "C08"
$ws3.Range("E11").Value = "Nu exista mesaje de eroare la aparitia erorilor"
[void]$ws3.Range("E11").Select()

# --- Sheet: Tool-basedCodeAnalysis ---
$ws4 = $wb.Worksheets.Item("Tool-basedCodeAnalysis")
[void]$ws4.Activate()
[void]$ws4.Range("D12").Select()
